$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each tuple: row, new DAMSLTag (column I), new DialogAct (column J)
$changes = @(
    @(23,  "sv", "Statement-opinion"),
    @(32,  "sv", "Statement-opinion"),
    @(43,  "sd", "Statement-non-opinion"),
    @(59,  "%",  "Uninterpretable"),
    @(77,  "sd", "Statement-non-opinion"),
    @(82,  "sv", "Statement-opinion"),
    @(85,  "aa", "Agree/Accept"),
    @(87,  "sd", "Statement-non-opinion"),
    @(91,  "sd", "Statement-non-opinion"),
    @(100, "sd", "Statement-non-opinion"),
    @(103, "qy", "Yes-No-Question"),
    @(108, "sv", "Statement-opinion"),
    @(111, "sv", "Statement-opinion"),
    @(112, "sv", "Statement-opinion"),
    @(117, "sv", "Statement-opinion")
)

foreach ($change in $changes) {
    $row = $change[0]
    $tag = $change[1]
    $act = $change[2]
    $ws.Cells.Item($row, 9).Value = $tag
    $ws.Cells.Item($row, 10).Value = $act
}
